$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("External driver companies")

# Insert a new column before column C (shifts National/International shift
# counts and the Driver starting address column one place to the right).
$ws.Columns.Item(3).Insert()

# Work around an engine quirk where inserting a column leaves a phantom
# column beyond the worksheet's last column (16384), which would otherwise
# get serialized as an out-of-range <col max="16385"/>. Deleting the true
# last column collapses that phantom back out.
$ws.Columns.Item(16384).Delete()

# New header for the inserted column.
$ws.Cells.Item(1, 3).Value = "Allows hotel stays?"

# Set the boolean "Allows hotel stays?" values for each company.
$ws.Cells.Item(2, 3).Value = $false
$ws.Cells.Item(3, 3).Value = $true
$ws.Cells.Item(4, 3).Value = $true
$ws.Cells.Item(5, 3).Value = $false

# Update the active selection to match the saved state of the workbook.
$ws.Range("C6").Select()
